$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6559
$ws.Range("J3").Value = 6944
$ws.Range("I4").Value = 1774
$ws.Range("J4").Value = 1509
$ws.Range("J5").Value = 538
$ws.Range("J6").Value = 9218
$ws.Range("I7").Value = 26231
$ws.Range("J7").Value = 24768

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 427
$ws.Range("J3").Value = 467
$ws.Range("J6").Value = 544
$ws.Range("J7").Value = 1559

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 144
$ws.Range("J3").Value = 184
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 493

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 262
$ws.Range("J3").Value = 370
$ws.Range("J4").Value = 50
$ws.Range("J6").Value = 392
$ws.Range("J7").Value = 1120

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 358

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 222
$ws.Range("J3").Value = 256
$ws.Range("J5").Value = 30
$ws.Range("J7").Value = 757

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 181
$ws.Range("J3").Value = 173
$ws.Range("J6").Value = 220
$ws.Range("J7").Value = 617

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 201
$ws.Range("J7").Value = 722
$ws.Range("J8").Value = 1559
$ws.Range("J11").Value = 422
$ws.Range("J16").Value = 101
$ws.Range("J19").Value = 723
$ws.Range("J20").Value = 519
$ws.Range("J29").Value = 1351
$ws.Range("J31").Value = 242
$ws.Range("J33").Value = 1120
$ws.Range("J36").Value = 335
$ws.Range("J37").Value = 757
$ws.Range("J42").Value = 1063
$ws.Range("J43").Value = 213
$ws.Range("J48").Value = 280
$ws.Range("J49").Value = 158
$ws.Range("J50").Value = 150
$ws.Range("J52").Value = 624
$ws.Range("I55").Value = 314
$ws.Range("J55").Value = 379
$ws.Range("J57").Value = 111
$ws.Range("J63").Value = 85
$ws.Range("J64").Value = 165
$ws.Range("J65").Value = 617
$ws.Range("J67").Value = 939
$ws.Range("J71").Value = 81
$ws.Range("J73").Value = 237
$ws.Range("J76").Value = 366
$ws.Range("J77").Value = 179
$ws.Range("J78").Value = 295
$ws.Range("J79").Value = 695
$ws.Range("J83").Value = 493
$ws.Range("J85").Value = 1025
$ws.Range("J87").Value = 82
$ws.Range("J88").Value = 255
$ws.Range("J89").Value = 318
$ws.Range("J91").Value = 283
$ws.Range("J95").Value = 358
$ws.Range("J97").Value = 221
$ws.Range("J100").Value = 46
$ws.Range("I101").Value = 26231
$ws.Range("J101").Value = 24768

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J4").Value = 13
$ws.Range("J7").Value = 242

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 239
$ws.Range("J3").Value = 349
$ws.Range("J7").Value = 939

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 158

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 414
$ws.Range("J3").Value = 473
$ws.Range("J6").Value = 341
$ws.Range("J7").Value = 1351

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J4").Value = 43
$ws.Range("J6").Value = 137
$ws.Range("J7").Value = 280

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 173
$ws.Range("J7").Value = 723

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J6").Value = 199
$ws.Range("J7").Value = 366

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 226
$ws.Range("J6").Value = 567
$ws.Range("J7").Value = 1063

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 92
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 295

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 73
$ws.Range("J3").Value = 75
$ws.Range("I4").Value = 14
$ws.Range("J6").Value = 212
$ws.Range("I7").Value = 314
$ws.Range("J7").Value = 379

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J3").Value = 73
$ws.Range("J6").Value = 97

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J5").Value = 11
$ws.Range("J7").Value = 283

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 192
$ws.Range("J6").Value = 206
$ws.Range("J7").Value = 695

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J4").Value = 16
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J6").Value = 146
$ws.Range("J7").Value = 519

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 102
$ws.Range("J7").Value = 335

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 220
$ws.Range("J7").Value = 722

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J3").Value = 38
$ws.Range("J4").Value = 23
$ws.Range("J7").Value = 150

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 122
$ws.Range("J7").Value = 422

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 77
$ws.Range("J7").Value = 237

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J2").Value = 59
$ws.Range("J6").Value = 75
$ws.Range("J7").Value = 201

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 63
$ws.Range("J6").Value = 126
$ws.Range("J7").Value = 255

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 91
$ws.Range("J7").Value = 318

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J2").Value = 28
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 111

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 213

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 274
$ws.Range("J3").Value = 362
$ws.Range("J4").Value = 67
$ws.Range("J6").Value = 297
$ws.Range("J7").Value = 1025

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J3").Value = 59
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 182
$ws.Range("J6").Value = 266
$ws.Range("J7").Value = 624

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 101
